$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibitions)
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All types)

# Update "想去人数" (want-to-go count) column F on sheet 1 (展览)
$ws1.Range("F5").Value = 1811
$ws1.Range("F6").Value = 60
$ws1.Range("F7").Value = 16
$ws1.Range("F8").Value = 160
$ws1.Range("F9").Value = 2281
$ws1.Range("F10").Value = 113
$ws1.Range("F11").Value = 63
$ws1.Range("F13").Value = 1398
$ws1.Range("F14").Value = 494
$ws1.Range("F15").Value = 31
$ws1.Range("F16").Value = 307
$ws1.Range("F17").Value = 214
$ws1.Range("F18").Value = 13
$ws1.Range("F19").Value = 26
$ws1.Range("F21").Value = 58
$ws1.Range("F23").Value = 1
$ws1.Range("F24").Value = 69
$ws1.Range("F26").Value = 1419
$ws1.Range("F28").Value = 365
$ws1.Range("F29").Value = 121
$ws1.Range("F32").Value = 354

# Update "想去人数" (want-to-go count) column F on sheet 4 (全部类型)
$ws4.Range("F5").Value = 1811
$ws4.Range("F7").Value = 60
$ws4.Range("F8").Value = 16
$ws4.Range("F9").Value = 160
$ws4.Range("F10").Value = 2281
$ws4.Range("F11").Value = 113
$ws4.Range("F12").Value = 63
$ws4.Range("F14").Value = 1398
$ws4.Range("F15").Value = 494
$ws4.Range("F16").Value = 31
$ws4.Range("F17").Value = 307
$ws4.Range("F18").Value = 214
$ws4.Range("F19").Value = 13
$ws4.Range("F20").Value = 26
$ws4.Range("F22").Value = 58
$ws4.Range("F24").Value = 1
$ws4.Range("F25").Value = 69
$ws4.Range("F27").Value = 1419
$ws4.Range("F29").Value = 365
$ws4.Range("F30").Value = 121
$ws4.Range("F33").Value = 354

